$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row: "<Name>_old" -> "<Name>_FV2404" (columns A-J / 1-10)
#    and "<Name>_new" -> "<Name>_FV2410" (columns L-U / 12-21).
$bases = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $bases.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($bases[$i])_FV2404"
}
for ($i = 0; $i -lt $bases.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($bases[$i])_FV2410"
}

# 2) Turn the used range into an Excel Table ("Table1") so the header row
#    gets autofilter buttons and structured references.
$rng = $ws.Range("A1:U88")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes, $null)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze the header row (split/freeze above row 2).
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
